# "Todos los botones funcionales"
# Update the Usuarios sheet data: the "test/testeo/ttt/gggg" placeholder row
# becomes a new "qweqwe/qwe/qwe/qwe" test row, the fff login for the fgh user
# becomes qqq, the SebastianJerez row keeps its data but its password column
# (2724) is now entered as a real number instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (fgh user): username fff -> qqq
$ws.Range("B2").Value = "qqq"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1

# Row 3: replace the old test/testeo/ttt/gggg placeholder data with the new
# qweqwe/qwe/qwe/qwe test entry
$ws.Range("A3").Value = "qweqwe"
$ws.Range("B3").Value = "qwe"
$ws.Range("C3").Value = "qwe"
$ws.Range("D3").Value = "qwe"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2

# Row 4: SebastianJerez row now holds the password (2724) as a real number
$ws.Range("A4").Value = "SebastianJerez"
$ws.Range("B4").Value = "sebas"
$ws.Range("C4").Value = 2724
$ws.Range("D4").Value = "sebastianjs99@hotmail.com"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 3

# Best-fit-ish column widths for the name / user / email columns
# (ColumnWidth is pixel-quantized by the host, so these are the inputs
# that land closest to the target stored widths of 16.21875 / 24.21875)
$ws.Columns.Item(1).ColumnWidth = 15.3
$ws.Columns.Item(2).ColumnWidth = 23.3
$ws.Columns.Item(4).ColumnWidth = 23.3

# Move/restore the active selection like the author left it
$ws.Range("F7").Select()
